# Sprint 7 backlog update:
#  - Insert a new backlog row (row 4) for the "cactus immunity" item.
#  - Fill in the new "cactus reveals hand" item that already had a
#    reserved/blank row template (row 7 after the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above the current row 4 ("Tavolo di gioco" / ...),
# pushing the rest of the backlog down by one row.
$ws.Rows.Item(4).Insert()

# New backlog item: immunity after calling cactus.
$ws.Range("B4").Value = "Giocatore immune agli effetti dopo aver chiamato cactus"
$ws.Range("C4").Value = "Artegiani"
$ws.Range("D4").Value = 2

# Copy the "Initial Size Estimate" marker formatting (col G) from the row
# below, which used to be the original row 4, so the new row matches it.
$ws.Range("G5").Copy($ws.Range("G4"))

# The row that used to be the blank template row 6 is now row 7 - fill it
# in with the new "cactus reveals hand" backlog item.
$ws.Range("B7").Value = "Quando un giocatore chiama cactus la sua mano viene scoperta"
$ws.Range("C7").Value = "Artegiani"
$ws.Range("D7").Value = 2
